$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round coordinate values to whole numbers
$ws.Range("Q2").Value = 357963
$ws.Range("R2").Value = 6875239

# Clear the Starttid (Z2) and Sluttid (AB2) cells entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
